$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the existing row 170. This pushes the
# current rows 170-264 down to 172-266 (and the sheet dimension grows from
# A1:R264 to A1:R266), mirroring two new daily price records being added
# at the top of this date-sorted block.
$ws.Rows("170:171").Insert()

# --- New row 170: new daily record ---
$ws.Range("A170").Value = 3
$ws.Range("B170").Value = "Femacal de La Calera"
$ws.Range("C170").Value = "Coquimbo"
$ws.Range("D170").Value = 44518
$ws.Range("E170").Value = 5
$ws.Range("F170").Value = 100112032
$ws.Range("G170").Value = "Zapallo italiano"
$ws.Range("H170").Value = "Sin especificar"
$ws.Range("I170").Value = "Primera"
$ws.Range("J170").Value = 140
$ws.Range("K170").Value = 4000
$ws.Range("L170").Value = 4000
$ws.Range("M170").Value = 4000
$ws.Range("N170").Value = "$/caja 36 unidades"
$ws.Range("O170").Value = "Provincia de Quillota"
$ws.Range("P170").Value = 111
$ws.Range("Q170").Value = 36
$ws.Range("R170").Value = "Hortaliza"

# --- New row 171: new daily record ---
$ws.Range("A171").Value = 3
$ws.Range("B171").Value = "Femacal de La Calera"
$ws.Range("C171").Value = "Coquimbo"
$ws.Range("D171").Value = 44518
$ws.Range("E171").Value = 5
$ws.Range("F171").Value = 100112032
$ws.Range("G171").Value = "Zapallo italiano"
$ws.Range("H171").Value = "Sin especificar"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 115
$ws.Range("K171").Value = 7500
$ws.Range("L171").Value = 8000
$ws.Range("M171").Value = 7739
$ws.Range("N171").Value = "$/caja 70 unidades"
$ws.Range("O171").Value = "Región de Arica y Parinacota"
$ws.Range("P171").Value = 111
$ws.Range("Q171").Value = 70
$ws.Range("R171").Value = "Hortaliza"
